$d = $word.ActiveDocument

# 1. Heading "3.2 Train the Model" -> "3.2 Train the model"
#    (the leading "T" is already a separate run and is left untouched;
#    only the "rain the Model" tail - spanning the "rain the "/"M"/"odel"
#    runs - is matched and merged into a single run.)
$d.Content.Find.Execute("rain the Model", $true, $false, $false, $false, $false, $true, 1, $false, "rain the model", 2) | Out-Null

# 2. Fix the misspelled TOC entry "C. Loss Fuction" -> "C. Loss Function"
$d.Content.Find.Execute("Fuction", $true, $false, $false, $false, $false, $true, 1, $false, "Function", 2) | Out-Null

# 3. Remove the stray space before the comma in " , test.py" -> ", test.py"
$d.Content.Find.Execute(" , test.py", $true, $false, $false, $false, $false, $true, 1, $false, ", test.py", 2) | Out-Null
